# Insert a new row of weekly price data at row 258 in the "Poroto verde"
# (green bean) sheet. All existing rows from 258 downward shift to 259+.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 258, pushing rows 258:277 down to 259:278
$ws.Rows.Item(258).Insert()

# Populate the new row 258 with the latest weekly record
$ws.Range("A258").Value = 8
$ws.Range("B258").Value = "Terminal La Palmera de La Serena"
$ws.Range("C258").Value = "Coquimbo"
$ws.Range("D258").Value = 44826
$ws.Range("E258").Value = 4
$ws.Range("F258").Value = 100112031
$ws.Range("G258").Value = "Poroto verde"
$ws.Range("H258").Value = "Magnum"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 480
$ws.Range("K258").Value = 36000
$ws.Range("L258").Value = 37000
$ws.Range("M258").Value = 36500
$ws.Range("N258").Value = "$/malla 25 kilos"
$ws.Range("O258").Value = "Perú"
$ws.Range("P258").Value = 1460
$ws.Range("Q258").Value = 25
$ws.Range("R258").Value = "Hortaliza"
